$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.284.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.19%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.507.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4 - TetherUSD
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.25%  "

# Row 8 - USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.62%  "

# Row 10 - Avalanche
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.82%  "

# Row 11 - Chainlink
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +8.77%  "

# Row 12 - Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0817"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13 - TRON
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.37%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.900.89"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.509.09"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17 - Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.130.61"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.00%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.14"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.38%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.76"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.00%  "

# Row 21 - ShibaInu
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.91%  "

# Row 22 - ImmutableX
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.22%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.21%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "277.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.54%  "

# Row 25 - PancakeSwap
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.59%  "

# Row 26 - Dai
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.45%  "

# Row 28 - Toncoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.36%  "

# Row 29 - Cosmos
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.83"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.23%  "

# Rows 30-31 - InjectiveProtocol and Kaspa swap positions
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.52%  "

# Row 32 - OKB
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.10"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.51%  "

# Row 33 - Celestia
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.59"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.51%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.36"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.75%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.02%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.57%  "

# Row 37 - ARBITRUM
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.10%  "

# Row 38 - RenderToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.66"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.47%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40 - Stellar
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41 - Monero
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.45"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.27%  "

# Row 42 - WEMIXToken
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.65%  "

# Row 43 - EnergySwap
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.51"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.85%  "

# Row 44 - VeChain
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.22%  "

# Row 45 - Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.002.62"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.21%  "

# Row 46 - NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.46%  "

# Row 47 - Stacks
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.76%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.02%  "

# Row 49 - FraxShare
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.98"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.62%  "

# Row 50 - THORChain
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.96%  "

# Row 51 - BitcoinSV
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.10"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.38%  "
